$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6412155.5
$ws.Range("J17").Value = 6946330
$ws.Range("L17").Value = 20838990
$ws.Range("N17").Value = -20839326
$ws.Range("H41").Value = 664.125
$ws.Range("I41").Value = 211.41667
$ws.Range("J41").Value = 2022.25
$ws.Range("K41").Value = 211.41667
$ws.Range("L41").Value = 2022.25
$ws.Range("M41").Value = 228.58333
$ws.Range("N41").Value = -2902.25
$ws.Range("H51").Value = 38470244
$ws.Range("I51").Value = 50008708
$ws.Range("J51").Value = 8698.666999999999
$ws.Range("K51").Value = 50008708
$ws.Range("L51").Value = 8698.666999999999
$ws.Range("M51").Value = -50008224
$ws.Range("N51").Value = -9666.666999999999
$ws.Range("H98").Value = 1823.8679
$ws.Range("I98").Value = 1851.2
$ws.Range("K98").Value = 1851.2
$ws.Range("M98").Value = -353.2
$ws.Range("H122").Value = 1823.8679
$ws.Range("I122").Value = 1851.2
$ws.Range("K122").Value = 5553.6
$ws.Range("M122").Value = -3103.6
$ws.Range("H138").Value = 297196.38
$ws.Range("J138").Value = 465972.25
$ws.Range("L138").Value = 1397916.75
$ws.Range("N138").Value = -1408196.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3243.724
$ws.Range("I32").Value = 3243.724
$ws.Range("K32").Value = 3243.724
$ws.Range("M32").Value = -2956.724
$ws.Range("H61").Value = 2539.0334
$ws.Range("I61").Value = 1548.1666
$ws.Range("K61").Value = 1548.1666
$ws.Range("M61").Value = -1336.1666
$ws.Range("H74").Value = 127560.98
$ws.Range("I74").Value = 169383.1
$ws.Range("J74").Value = 2094.6365
$ws.Range("K74").Value = 169383.1
$ws.Range("L74").Value = 2094.6365
$ws.Range("M74").Value = -168509.1
$ws.Range("N74").Value = -3842.6365
$ws.Range("H77").Value = 127560.98
$ws.Range("I77").Value = 169383.1
$ws.Range("J77").Value = 2094.6365
$ws.Range("K77").Value = 846915.5
$ws.Range("L77").Value = 10473.1825
$ws.Range("M77").Value = -842547.5
$ws.Range("N77").Value = -19209.1825
$ws.Range("H110").Value = 2915.0833
$ws.Range("I110").Value = 2069.9033
$ws.Range("J110").Value = 8155.2
$ws.Range("K110").Value = 2069.9033
$ws.Range("L110").Value = 8155.2
$ws.Range("M110").Value = -24.90329999999994
$ws.Range("N110").Value = -12245.2
$ws.Range("H122").Value = 4950
$ws.Range("I122").Value = 4645.4707
$ws.Range("K122").Value = 13936.4121
$ws.Range("M122").Value = -11486.4121
$ws.Range("H132").Value = 1716.9464
$ws.Range("I132").Value = 1305
$ws.Range("J132").Value = 3868.2222
$ws.Range("K132").Value = 3915
$ws.Range("L132").Value = 11604.6666
$ws.Range("M132").Value = -1385
$ws.Range("N132").Value = -16664.6666
$ws.Range("H136").Value = 2539.0334
$ws.Range("I136").Value = 1548.1666
$ws.Range("K136").Value = 4644.4998
$ws.Range("M136").Value = -2094.4998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3587.7917
$ws.Range("I99").Value = 3155.4
$ws.Range("J99").Value = 5749.75
$ws.Range("K99").Value = 3155.4
$ws.Range("L99").Value = 5749.75
$ws.Range("M99").Value = -1657.4
$ws.Range("N99").Value = -8745.75
$ws.Range("H105").Value = 20003258
$ws.Range("I105").Value = 1252831.8
$ws.Range("K105").Value = 1252831.8
$ws.Range("M105").Value = -1251084.8
$ws.Range("H134").Value = 1329.579
$ws.Range("I134").Value = 952.4194
$ws.Range("J134").Value = 2999.8572
$ws.Range("K134").Value = 2857.2582
$ws.Range("L134").Value = 8999.571599999999
$ws.Range("M134").Value = -322.2582000000002
$ws.Range("N134").Value = -14069.5716

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2426.359
$ws.Range("I58").Value = 1673.5555
$ws.Range("K58").Value = 1673.5555
$ws.Range("M58").Value = -1470.5555
$ws.Range("H94").Value = 2061.1765
$ws.Range("I94").Value = 1766.5
$ws.Range("J94").Value = 2221.9092
$ws.Range("K94").Value = 1766.5
$ws.Range("L94").Value = 2221.9092
$ws.Range("M94").Value = -1315.5
$ws.Range("N94").Value = -3123.9092
$ws.Range("H99").Value = 8382.388999999999
$ws.Range("I99").Value = 8759.066000000001
$ws.Range("J99").Value = 6499
$ws.Range("K99").Value = 8759.066000000001
$ws.Range("L99").Value = 6499
$ws.Range("M99").Value = -7261.066000000001
$ws.Range("N99").Value = -9495
$ws.Range("H126").Value = 8382.388999999999
$ws.Range("I126").Value = 8759.066000000001
$ws.Range("J126").Value = 6499
$ws.Range("K126").Value = 26277.198
$ws.Range("L126").Value = 19497
$ws.Range("M126").Value = -23807.198
$ws.Range("N126").Value = -24437
$ws.Range("H132").Value = 3047.0356
$ws.Range("I132").Value = 2392.3157
$ws.Range("K132").Value = 7176.9471
$ws.Range("M132").Value = -4646.9471
$ws.Range("H136").Value = 2426.359
$ws.Range("I136").Value = 1673.5555
$ws.Range("K136").Value = 5020.666499999999
$ws.Range("M136").Value = -2470.666499999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 9464
$ws.Range("J131").Value = 2391.9
$ws.Range("L131").Value = 7175.700000000001
$ws.Range("N131").Value = -17255.7
$ws.Range("H134").Value = 5774.1
$ws.Range("I134").Value = 2217.625
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 6652.875
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -1582.875
$ws.Range("N134").Value = -70140
$ws.Range("H139").Value = 2603.739
$ws.Range("I139").Value = 1834.7142
$ws.Range("K139").Value = 5504.142599999999
$ws.Range("M139").Value = -364.1425999999992
$ws.Range("H141").Value = 9948.5
$ws.Range("I141").Value = 9948.5
$ws.Range("K141").Value = 29845.5
$ws.Range("M141").Value = -24665.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 768.7586
$ws.Range("I97").Value = 428.125
$ws.Range("J97").Value = 2403.8
$ws.Range("K97").Value = 428.125
$ws.Range("L97").Value = 2403.8
$ws.Range("M97").Value = 67.875
$ws.Range("N97").Value = -3395.8
$ws.Range("H126").Value = 8155.615
$ws.Range("I126").Value = 6821.1
$ws.Range("J126").Value = 12604
$ws.Range("K126").Value = 20463.3
$ws.Range("L126").Value = 37812
$ws.Range("M126").Value = -17993.3
$ws.Range("N126").Value = -42752
$ws.Range("H132").Value = 1397.5883
$ws.Range("I132").Value = 1397.5883
$ws.Range("K132").Value = 4192.7649
$ws.Range("M132").Value = -1662.7649
$ws.Range("H141").Value = 57398
$ws.Range("J141").Value = 57398
$ws.Range("L141").Value = 57398
$ws.Range("N141").Value = -67758

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5260.5806
$ws.Range("I132").Value = 4022.0386
$ws.Range("K132").Value = 12066.1158
$ws.Range("M132").Value = -9536.1158
$ws.Range("H136").Value = 3193.6086
$ws.Range("I136").Value = 2886.6572
$ws.Range("K136").Value = 8659.971600000001
$ws.Range("M136").Value = -6109.971600000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 626.15
$ws.Range("I113").Value = 572.6923
$ws.Range("K113").Value = 1718.0769
$ws.Range("M113").Value = 451.9231
